$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, shifting the existing rows 118:257 down to 119:258.
$ws.Rows(118).Insert()

# Populate the newly inserted row 118 with the new weekly data point.
$ws.Cells.Item(118, 1).Value = 8
$ws.Cells.Item(118, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44664
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = 100112012
$ws.Cells.Item(118, 7).Value = "Espinaca"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 2520
$ws.Cells.Item(118, 11).Value = 550
$ws.Cells.Item(118, 12).Value = 600
$ws.Cells.Item(118, 13).Value = 575
$ws.Cells.Item(118, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(118, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(118, 16).Value = 1150
$ws.Cells.Item(118, 17).Value = 0.5
$ws.Cells.Item(118, 18).Value = "Hortaliza"
